# Regenerate orders with updated distance/size codes.
#
# The experiment's distance and size condition labels were renumbered:
#   D80 -> D86
#   D64 -> D69
#   D51 -> D55
#   S30 -> S31
# (S20/S25 are unchanged.)
#
# These substrings show up embedded inside many text values across the
# sheet (condition names, stimulus filenames, and the Distance/Size lookup
# columns), so the fix is applied as a substring replacement on every
# text cell in the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val.Replace("D80", "D86").Replace("D64", "D69").Replace("D51", "D55").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
